$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The K column holds a date/time timestamp for each transaction row (rows 2-32).
# All of them move forward by exactly one day (2024-05-29 -> 2024-05-30),
# i.e. the underlying serial value goes from 45441.615127314813 to 45442.615127314813.
$ws.Range("K2:K32").Value = 45442.615127314813
